$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 data
$ws.Cells.Item(5, 1).Value = 42607.887337962966

$ws.Cells.Item(5, 2).Value = 12
$ws.Cells.Item(5, 3).Value = 54
$ws.Cells.Item(5, 4).Value = 45
$ws.Cells.Item(5, 5).Value = 66
$ws.Cells.Item(5, 6).Value = 33
$ws.Cells.Item(5, 7).Value = 14367
$ws.Cells.Item(5, 8).Value = 14986
$ws.Cells.Item(5, 9).Value = 2336
$ws.Cells.Item(5, 10).Value = 271
$ws.Cells.Item(5, 11).Value = 227
$ws.Cells.Item(5, 12).Value = 10
$ws.Cells.Item(5, 13).Value = 5
$ws.Cells.Item(5, 14).Value = "Noun"
